$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.77102397423476254
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

$ws.Range("B3").Value = -0.90133221548553433
$ws.Range("C3").Value = 0.37983388905294763
$ws.Range("D3").Value = -0.79127357551109223
$ws.Range("E3").Value = 2.2441822021855384

$ws.Range("B1:E3").Select()
